$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 235, shifting the existing rows 235-329 down to 236-330
$ws.Rows.Item(235).Insert()

# Populate the newly inserted row with the Aracatuba, Brazil colo data
$ws.Cells.Item(235, 1).Value = "ARU"
$ws.Cells.Item(235, 2).Value = "Aracatuba, Brazil"
$ws.Cells.Item(235, 3).Value = -21.1413002014
$ws.Cells.Item(235, 4).Value = -50.4247016907
$ws.Cells.Item(235, 5).Value = "BR"
$ws.Cells.Item(235, 6).Value = "South America"
$ws.Cells.Item(235, 7).Value = "Aracatuba"

# Normalize the "colo" column cell formatting so it matches the bold/bordered
# style used by every other row in column A.
$ws.Cells.Item(235, 1).Borders.LineStyle = 1
